$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Large multi-line text blocks (here-strings) for the new columns
$imageSW1 = @"
Cisco IOS Software, vios_l2 Software (vios_l2-ADVENTERPRISEK9-M), Version 15.2(CML_NIGHTLY_20190423)FLO_DSGS7, EARLY DEPLOYMENT DEVELOPMENT BUILD, synced to  V152_6_0_81_E
Technical Support: http://www.cisco.com/techsupport
Copyright (c) 1986-2019 by Cisco Systems, Inc.
Compiled Tue 23-Apr-19 04:48 by mmen
ROM: Bootstrap program is IOSv
SW1 uptime is 6 hours, 41 minutes
System returned to ROM by reload
System image file is "flash0:/vios_l2-adventerprisek9-m"
Last reload reason: Unknown reason
This product contains cryptographic features and is subject to United
States and local country laws governing import, export, transfer and
use. Delivery of Cisco cryptographic products does not imply
third-party authority to import, export, distribute or use encryption.
Importers, exporters, distributors and users are responsible for
compliance with U.S. and local country laws. By using this product you
agree to comply with applicable laws and regulations. If you are unable
to comply with U.S. and local laws, return this product immediately.
A summary of U.S. laws governing Cisco cryptographic products may be found at:
http://www.cisco.com/wwl/export/crypto/tool/stqrg.html
If you require further assistance please contact us by sending email to
export@cisco.com.
Cisco IOSv () processor (revision 1.0) with 935161K/111616K bytes of memory.
Processor board ID 9BQHFDLPE7P
1 Virtual Ethernet interface
8 Gigabit Ethernet interfaces
DRAM configuration is 72 bits wide with parity disabled.
256K bytes of non-volatile configuration memory.
2097152K bytes of ATA System CompactFlash 0 (Read/Write)
0K bytes of ATA CompactFlash 1 (Read/Write)
0K bytes of ATA CompactFlash 2 (Read/Write)
0K bytes of ATA CompactFlash 3 (Read/Write)
Configuration register is 0x101

"@

$imageSW2 = @"
Cisco IOS Software, vios_l2 Software (vios_l2-ADVENTERPRISEK9-M), Version 15.2(CML_NIGHTLY_20190423)FLO_DSGS7, EARLY DEPLOYMENT DEVELOPMENT BUILD, synced to  V152_6_0_81_E
Technical Support: http://www.cisco.com/techsupport
Copyright (c) 1986-2019 by Cisco Systems, Inc.
Compiled Tue 23-Apr-19 04:48 by mmen
ROM: Bootstrap program is IOSv
SW2 uptime is 6 hours, 25 minutes
System returned to ROM by reload
System image file is "flash0:/vios_l2-adventerprisek9-m"
Last reload reason: Unknown reason
This product contains cryptographic features and is subject to United
States and local country laws governing import, export, transfer and
use. Delivery of Cisco cryptographic products does not imply
third-party authority to import, export, distribute or use encryption.
Importers, exporters, distributors and users are responsible for
compliance with U.S. and local country laws. By using this product you
agree to comply with applicable laws and regulations. If you are unable
to comply with U.S. and local laws, return this product immediately.
A summary of U.S. laws governing Cisco cryptographic products may be found at:
http://www.cisco.com/wwl/export/crypto/tool/stqrg.html
If you require further assistance please contact us by sending email to
export@cisco.com.
Cisco IOSv () processor (revision 1.0) with 935161K/111616K bytes of memory.
Processor board ID 9HR276JJI1R
1 Virtual Ethernet interface
8 Gigabit Ethernet interfaces
DRAM configuration is 72 bits wide with parity disabled.
256K bytes of non-volatile configuration memory.
2097152K bytes of ATA System CompactFlash 0 (Read/Write)
0K bytes of ATA CompactFlash 1 (Read/Write)
0K bytes of ATA CompactFlash 2 (Read/Write)
0K bytes of ATA CompactFlash 3 (Read/Write)
Configuration register is 0x101

"@

$cdpSW1 = @"
Capability Codes: R - Router, T - Trans Bridge, B - Source Route Bridge
                  S - Switch, H - Host, I - IGMP, r - Repeater, P - Phone, 
                  D - Remote, C - CVTA, M - Two-port Mac Relay 
Device ID        Local Intrfce     Holdtme    Capability  Platform  Port ID
MikroTik         Gig 0/1           100               R    MikroTik  LAN/ether2
SW2.home.net     Gig 0/1           143             R S I            Gig 0/1
Total cdp entries displayed : 2
"@

$cdpSW2 = @"
Capability Codes: R - Router, T - Trans Bridge, B - Source Route Bridge
                  S - Switch, H - Host, I - IGMP, r - Repeater, P - Phone, 
                  D - Remote, C - CVTA, M - Two-port Mac Relay 
Device ID        Local Intrfce     Holdtme    Capability  Platform  Port ID
MikroTik         Gig 0/1           100               R    MikroTik  LAN/ether2
SW1.home.net     Gig 0/1           125             R S I            Gig 0/1
Total cdp entries displayed : 2
"@

$intDescr = @"
Interface                      Status         Protocol Description
Gi0/0                          up             up       TEST2022
Gi0/1                          up             up       TEST2022
Gi0/2                          up             up       TEST2022
Gi0/3                          up             up       TEST2022
Gi1/0                          up             up       TEST2022
Gi1/1                          up             up       TEST2022
Gi1/2                          up             up       TEST2022
Gi1/3                          up             up       TEST2022
Vl10                           up             up       TEST2022
"@

# Row 2 - SW1
$ws.Range("E2").Value = "SW1"
$ws.Range("F2").Value = "1.1.1.190"
$ws.Range("H2").Value = "cisco_ios"
$ws.Range("I2").Value = "ssh"
$ws.Range("J2").Value = "hosts_if.xlsx"
$ws.Range("K2").Value = "interface_cmd.template"
$ws.Range("L2").Value = "no errors"
$ws.Range("M2").Value = "no error"
$ws.Range("O2").Value = "vios_l2-adventerprisek9-m"
$ws.Range("P2").Value = $imageSW1
$ws.Range("W2").Value = $cdpSW1
$ws.Range("X2").Value = "% LLDP is not enabled"
$ws.Range("Y2").Value = $intDescr

# Row 3 - SW2
$ws.Range("E3").Value = "SW2"
$ws.Range("F3").Value = "1.1.1.195"
$ws.Range("H3").Value = "cisco_ios"
$ws.Range("I3").Value = "ssh"
$ws.Range("J3").Value = "hosts_if.xlsx"
$ws.Range("K3").Value = "interface_cmd.template"
$ws.Range("L3").Value = "no errors"
$ws.Range("M3").Value = "no error"
$ws.Range("O3").Value = "vios_l2-adventerprisek9-m"
$ws.Range("P3").Value = $imageSW2
$ws.Range("W3").Value = $cdpSW2
$ws.Range("X3").Value = "% LLDP is not enabled"
$ws.Range("Y3").Value = $intDescr

# Reset auto row height (engine auto-sizes rows on multi-line values; the
# source workbook keeps default row height, so undo that side effect)
$ws.Rows.Item(2).AutoFit()
$ws.Rows.Item(3).AutoFit()

# Selection, to match saved view state
$ws.Range("K15").Select() | Out-Null
